# Actualización automática 2025-11-18 16:30:08
#
# 1) "VENTAS POR GRUPO": a new client row (OFICINA-CATAECSA / LINDAO RODRIGUEZ
#    ANTONIO COLON) is inserted at row 306, pushing every subsequent row down
#    by one; a handful of unrelated numeric cells are also corrected; the
#    trailing "x de N" summary row is refreshed for the new row count.
# 2) "VENTA MENSUAL": mirrors the same new client insertion (row 310) plus a
#    few unrelated numeric corrections and an updated totals row.
# 3) "CUMPLIMIENTO MENSUAL": recomputed rollup values that depend on the
#    numeric corrections made on the first two sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert the new client row; Excel shifts rows 306:352 down to 307:353 and
# clones the formatting of the row above (s="2" on the numeric columns).
$ws1.Rows.Item(306).Insert()

$ws1.Cells.Item(306, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(306, 2).Value = "LINDAO RODRIGUEZ ANTONIO COLON"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(306, $c).Value = 0
}

# Unrelated numeric corrections (rows are all above the insertion point, so
# their row numbers are unaffected by the shift).
$ws1.Cells.Item(80, 11).Value = 417.36    # K80
$ws1.Cells.Item(98, 13).Value = 4971.51   # M98
$ws1.Cells.Item(294, 8).Value = 423.9     # H294
$ws1.Cells.Item(294, 9).Value = 503.33    # I294

# Refresh the trailing "x de 350" -> "x de 351" summary row, now at row 353.
# Column K's count also increases (2 -> 3) because K80 became non-zero.
$ws1.Cells.Item(353, 3).Value = "3 de 351"
$ws1.Cells.Item(353, 4).Value = "14 de 351"
$ws1.Cells.Item(353, 5).Value = "6 de 351"
$ws1.Cells.Item(353, 6).Value = "0 de 351"
$ws1.Cells.Item(353, 7).Value = "0 de 351"
$ws1.Cells.Item(353, 8).Value = "5 de 351"
$ws1.Cells.Item(353, 9).Value = "10 de 351"
$ws1.Cells.Item(353, 10).Value = "0 de 351"
$ws1.Cells.Item(353, 11).Value = "3 de 351"
$ws1.Cells.Item(353, 12).Value = "22 de 351"
$ws1.Cells.Item(353, 13).Value = "43 de 351"
$ws1.Cells.Item(353, 14).Value = "2 de 351"
$ws1.Cells.Item(353, 15).Value = "0 de 351"
$ws1.Cells.Item(353, 16).Value = "1 de 351"
$ws1.Cells.Item(353, 17).Value = "0 de 351"
$ws1.Cells.Item(353, 18).Value = "0 de 351"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(310).Insert()

$ws2.Cells.Item(310, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(310, 2).Value = "LINDAO RODRIGUEZ ANTONIO COLON"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(310, $c).Value = 0
}

$ws2.Cells.Item(80, 6).Value = 860.8     # F80
$ws2.Cells.Item(98, 6).Value = 6698.96   # F98
$ws2.Cells.Item(298, 6).Value = 927.23   # F298

# Totals row, now at row 357; only the GRANITO ("F") column total changes,
# because only F-column source cells were corrected above.
$ws2.Cells.Item(357, 6).Value = 140204.23

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Cells.Item(22, 4).Value = 10.45
$ws3.Cells.Item(22, 5).Value = 9905.549999999999
$ws3.Cells.Item(22, 6).Value = 0.001053852359822509

$ws3.Cells.Item(36, 4).Value = 13497.27
$ws3.Cells.Item(36, 5).Value = 51446.73
$ws3.Cells.Item(36, 6).Value = 0.2078293606799704

$ws3.Cells.Item(63, 4).Value = 4518.7
$ws3.Cells.Item(63, 5).Value = -4518.7

$ws3.Cells.Item(77, 4).Value = 140069.77
$ws3.Cells.Item(77, 5).Value = 277178.9097415454
$ws3.Cells.Item(77, 6).Value = 0.3356985337539303

Write-Output "edit complete"
